# Update the dSF column (F) values for rows 2-10 and 12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = -7
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = -10
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = -2
